{"js": "// Revert the \"acknowledgement of service\" rename: change the document\n// title from \"Acknowledgement of claim\" back to \"Acknowledgement of service\".\n//\n// The title run originally reads \"Acknowledgement\" + \" of \" + \"claim\".\n// We only need to touch the \" of claim\" portion (turning it into\n// \" of service\"), leaving the leading \"Acknowledgement\" run untouched.\nconst body = context.document.body;\nconst results = body.search(\" of claim\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items,text\");\nawait context.sync();\n\nif (results.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for ' of claim', found \" + results.items.length);\n}\n\nresults.items[0].insertText(\" of service\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Revert the \"acknowledgement of service\" rename: change the document\n# title from \"Acknowledgement of claim\" back to \"Acknowledgement of service\".\n#\n# The title run originally reads \"Acknowledgement\" + \" of \" + \"claim\".\n# We only touch the trailing \"claim\" word (turning it into \"service\"),\n# leaving the \"Acknowledgement\" and \" of \" runs untouched.\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$find = $range.Find\n$find.ClearFormatting()\n$find.Text = \"claim\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $true\n\n$found = $find.Execute()\nif (-not $found) {\n    throw \"Could not find the word 'claim' in the document title.\"\n}\n\n$range.Text = \"service\"\n"}
